# ajustes 1ra quincena mar
$wb = $excel.ActiveWorkbook

# --- mar2025: enter payment amounts for 1st quincena of march ---
$wsMar = $wb.Worksheets.Item("mar2025")

$wsMar.Range("C2").Value = 32500
$wsMar.Range("C3").Value = 65000
$wsMar.Range("C4").Value = 65000
$wsMar.Range("C5").Value = 65000
$wsMar.Range("C7").Value = 65000
$wsMar.Range("C9").Value = 65000
$wsMar.Range("C10").Value = 65000
$wsMar.Range("D10").Value = 65000
$wsMar.Range("C14").Value = 65000
$wsMar.Range("D14").Value = 65000
$wsMar.Range("C15").Value = 65000
$wsMar.Range("C17").Value = 65000
$wsMar.Range("C20").Value = 65000
$wsMar.Range("C21").Value = 50000
$wsMar.Range("D21").Value = 50000
$wsMar.Range("C24").Value = 65000

# leave the cursor on C10 for this sheet
$wsMar.Range("C10").Select()

# --- feb2025: move the selection from D24 to D25 and keep it the active tab ---
$wsFeb = $wb.Worksheets.Item("feb2025")
$wsFeb.Range("D25").Select()
